# ATTENDANCE.xlsx edit script
# Rebuilds Sheet1 from the 5-row/15-col original into the expanded
# 8-row / 31-day attendance sheet (plus a stray "s" in V13), adds a
# "Phone Number" column, a "TOTAL" column, re-widens columns and flips
# the page to portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new column at C ("Phone Number") - this pushes the old
#    C (ROOM NO), D..N (dates) and N (Percentage) one column to the
#    right, preserving their values/styles automatically.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).Insert()

# New header for the inserted column
$ws.Range("C1").Value = "Phone Number"

# ---------------------------------------------------------------------
# 2. Extend the header row with the rest of March's dates (the sheet
#    previously only went to day 10 / column M [old] -> N [new]).
#    Columns E1:N1 (days 1-10) already hold 44986-44995 after the
#    column insert; append days 11-31 in O1:AI1, each using the same
#    date style (s=1) as the existing date header cells (E1 is one).
#    Overwriting O1 also clobbers the old "Percentage" label that the
#    column insert shifted into it - re-home it at AK1 afterwards.
# ---------------------------------------------------------------------
$ws.Range("E1").Copy() | Out-Null
$ws.Range("O1:AI1").PasteSpecial(-4122) | Out-Null

$day = 44996
for ($col = 15; $col -le 35; $col++) {
    $ws.Cells.Item(1, $col).Value = $day
    $day = $day + 1
}

# TOTAL header, and the re-homed "Percentage" header
$ws.Range("AJ1").Value = "TOTAL"
$ws.Range("AK1").Value = "Percentage"

# ---------------------------------------------------------------------
# 3. Phone numbers for the existing 5 students (column C)
# ---------------------------------------------------------------------
$ws.Range("C2").Value = 9174741646
$ws.Range("C3").Value = 7680011239
$ws.Range("C4").Value = 8604542123
$ws.Range("C5").Value = 7680011239
$ws.Range("C6").Value = 9284023616

# Row 2/3 used to hold a plain "1" in H2/H3 (old) / I2/I3 (new, post-
# insert) marking attendance - both get overwritten by "P" marks below
# (along with the rest of the attendance marks). P2 already carries
# the left-over date-style formatting (s=1), shifted in by the column
# insert from the old (empty) O2 cell - no extra styling needed there.
# ---------------------------------------------------------------------
# 4. Attendance ("P") marks for rows 2, 3, 7, 8
# ---------------------------------------------------------------------
$presentCells = @(
    "G2","I2","K2","L2","P2","Q2","S2",
    "E3","F3","G3","H3","I3","M3","O3","R3","S3",
    "S7",
    "G8","H8"
)
foreach ($addr in $presentCells) {
    $ws.Range($addr).Value = "P"
}

# ---------------------------------------------------------------------
# 5. TOTAL / Percentage columns (AJ/AK) for every student row
# ---------------------------------------------------------------------
$ws.Range("AJ2").Value = 7
$ws.Range("AK2").Value = 46.66666666666666

$ws.Range("AJ3").Value = 9
$ws.Range("AK3").Value = 60

$ws.Range("AJ4").Value = 0
$ws.Range("AK4").Value = 0

$ws.Range("AJ5").Value = 0
$ws.Range("AK5").Value = 0

$ws.Range("AJ6").Value = 1
$ws.Range("AK6").Value = 6.666666666666667

# ---------------------------------------------------------------------
# 6. Two brand-new students (rows 7 & 8)
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "RISHABH SINGH"
$ws.Range("B7").Value = 220953044
$ws.Range("C7").Value = 8604542123
$ws.Range("D7").Value = "F-12"
$ws.Range("AJ7").Value = 1
$ws.Range("AK7").Value = 6.666666666666667

$ws.Range("A8").Value = "ALLU ARJUN"
$ws.Range("B8").Value = 22095032
$ws.Range("C8").Value = 9174741646
$ws.Range("D8").Value = "G-12"
$ws.Range("AJ8").Value = 3
$ws.Range("AK8").Value = 20

# ---------------------------------------------------------------------
# 7. Stray "s" typed into V13
# ---------------------------------------------------------------------
$ws.Range("V13").Value = "s"

# ---------------------------------------------------------------------
# 8. Column widths
#    (ColumnWidth is quantized to whole pixels in this engine, so the
#    inputs below are the closest achievable match to the target
#    stored widths.)
# ---------------------------------------------------------------------
$ws.Range("A:A").ColumnWidth = 18.5
$ws.Range("B:C").ColumnWidth = 21.166666666666668
$ws.Range("D:D").ColumnWidth = 8.666666666666666
$ws.Range("E:N").ColumnWidth = 9.5
$ws.Range("O:O").ColumnWidth = 10.666666666666666
$ws.Range("P:AI").ColumnWidth = 9.5
$ws.Range("AK:AK").ColumnWidth = 14.166666666666666

# ---------------------------------------------------------------------
# 9. View: select V13, scroll so column V is left-most visible column
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 22
$win.ScrollRow = 1
$ws.Range("V13").Select() | Out-Null

# ---------------------------------------------------------------------
# 10. Page orientation -> portrait
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
